$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.542.01"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.322.75"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("Z1").Formula = "=`"576.35`""
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("Z1").Formula = "=`"173.95`""
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("Z1").Formula = "=`"0.999`""
$ws.Range("Z1").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "3.321.64"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("Z1").Formula = "=`"45.96`""
$ws.Range("Z1").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("Z1").Formula = "=`"703.28`""
$ws.Range("Z1").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "3.863.80"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "67.580.15"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "3.325.66"
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("Z1").Formula = "=`"5.36`""
$ws.Range("Z1").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E23").Value = "  +4.89%  "
$ws.Range("Z1").Formula = "=`"16.90`""
$ws.Range("Z1").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("Z1").Formula = "=`"98.53`""
$ws.Range("Z1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("Z1").Formula = "=`"3.86`""
$ws.Range("Z1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("Z1").Formula = "=`"9.35`""
$ws.Range("Z1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("Z1").Formula = "=`"33.18`""
$ws.Range("Z1").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("Z1").Formula = "=`"8.50`""
$ws.Range("Z1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("Z1").Formula = "=`"7.08`""
$ws.Range("Z1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E31").Value = "  +7.17%  "
$ws.Range("Z1").Formula = "=`"568.31`""
$ws.Range("Z1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("Z1").Formula = "=`"10.94`""
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("Z1").Formula = "=`"0.999`""
$ws.Range("Z1").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "3.696.25"
$ws.Range("E36").Value = "  -4.52%  "
$ws.Range("Z1").Formula = "=`"56.83`""
$ws.Range("Z1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("Z1").Formula = "=`"3.32`""
$ws.Range("Z1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  +6.45%  "
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("Z1").Formula = "=`"3.16`""
$ws.Range("Z1").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("Z1").Formula = "=`"2.62`""
$ws.Range("Z1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "0.0₃0670"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("Z1").Formula = "=`"0.334`""
$ws.Range("Z1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("Z1").Formula = "=`"3.25`""
$ws.Range("Z1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("Z1").Formula = "=`"0.0405`""
$ws.Range("Z1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("Z1").Formula = "=`"2.69`""
$ws.Range("Z1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E47").Value = "  +8.43%  "
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("E50").Value = "  -5.01%  "
$ws.Range("Z1").Formula = "=`"129.09`""
$ws.Range("Z1").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E51").Value = "  -0.01%  "
